# Computers in the Libraries.pptx -- apply commit:
#   - The "Data formats:" slide (originally slide index 4) moves to the
#     end of the deck.
#   - The "Correlations and hypotheses:" slide (originally slide index 5)
#     moves up to index 4.
#   - A new "Basic Tasks" slide is inserted at index 5 (taking over what
#     used to be the physical slide5 content-wise).
#   - Three brand-new slides are appended before the relocated
#     "Data formats:" slide: "Further Tasks", "Final tasks", "Assignments".
#
# Net visible slide order after edit:
#   1. Computers in the Libraries
#   2. Questions raised:
#   3. Data sets
#   4. Correlations and hypotheses:
#   5. Basic Tasks
#   6. Further Tasks
#   7. Final tasks
#   8. Assignments
#   9. Data formats:

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Step 1: slide 4 ("Data formats:") becomes "Correlations and hypotheses:"
# (it inherits exactly what used to live on slide 5).
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Correlations and hypotheses:"

$body4 = $s4.Shapes.Item(2).TextFrame
$tr4 = $body4.TextRange
$tr4.Text = "Weather: positively or negatively correlated with usage."
$null = $tr4.InsertAfter("`rStudents prefer to stay at home on days where it is cold/raining.")
$null = $tr4.InsertAfter("`rHours:  some machines are available 24/5. Others are not.")
$null = $tr4.InsertAfter("`rAn increase in availability necessarily indicates a commensurate increase in usage.")
$null = $tr4.InsertAfter("`rLocation: machines on group study v. quiet floors.")
$null = $tr4.InsertAfter("`rMachines on group study floors are used less frequently.")
$null = $tr4.InsertAfter("`rConfiguration: single or dual monitors.")
$null = $tr4.InsertAfter("`rMachines with dual monitors are preferred.")
$null = $tr4.InsertAfter("`rSchedules: exam times & breaks can change the data")
$null = $tr4.InsertAfter("`rUsage peaks during exam times, and drops off during breaks.")

# indent the "sub-bullet" paragraphs (2nd, 4th, 6th, 8th, 10th)
$tr4.Paragraphs(2).IndentLevel = 2
$tr4.Paragraphs(4).IndentLevel = 2
$tr4.Paragraphs(6).IndentLevel = 2
$tr4.Paragraphs(8).IndentLevel = 2
$tr4.Paragraphs(10).IndentLevel = 2

# shrink-text-on-overflow autofit (matches normAutofit on this shape)
$body4.AutoSize = 2

# ---------------------------------------------------------------------
# Step 2: slide 5 ("Correlations and hypotheses:") becomes "Basic Tasks"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Basic Tasks"

$body5 = $s5.Shapes.Item(2).TextFrame
$tr5 = $body5.TextRange
$tr5.Text = "1.a) Calculate total time computer in use per hour,      "
$null = $tr5.InsertAfter("`r1.b) `"`" per day,    ")
$null = $tr5.InsertAfter("`r1.c `"`" per week,    ")
$null = $tr5.InsertAfter("`r1.d `"`" per semester    ")
$null = $tr5.InsertAfter("`r2. Compare 1) against frequency of inclement rainfall   ")
$null = $tr5.InsertAfter("`r3. Compare 1) against times/days/weeks of exams    ")
$null = $tr5.InsertAfter("`r4. Compare 1) against ")

$p7 = $tr5.Paragraphs(7)
$null = $p7.InsertAfter("libary")
$null = $p7.InsertAfter(" ")
$null = $p7.InsertAfter("populaton")

$body5.AutoSize = 2

# ---------------------------------------------------------------------
# Step 3: new slide 6 -- "Further Tasks"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Further Tasks"

$body6 = $s6.Shapes.Item(2).TextFrame
$tr6 = $body6.TextRange
$tr6.Text = "1. Label each computer with terms 'window', 'dual monitor', 'quiet floor', or any combination of the three.    "
$null = $tr6.InsertAfter("`r2. Run all the basic tasks on each ")
$null = $tr6.InsertAfter("seperate")
$null = $tr6.InsertAfter(" category of computer.    ")
$null = $tr6.InsertAfter("`r3. Graph all the above data together? Draw conclusions: Which computers are used most often and when? Does weather affect how often/how much computers are used?    ")

$body6.AutoSize = 2

# ---------------------------------------------------------------------
# Step 4: new slide 7 -- "Final tasks"
# ---------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 2)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Final tasks"

$tr7 = $s7.Shapes.Item(2).TextFrame.TextRange
$tr7.Text = "1. Make a heat map of 'best' computers    "
$null = $tr7.InsertAfter("`r2. Machine Learning - predict which computers will be used when and how often")

# ---------------------------------------------------------------------
# Step 5: new slide 8 -- "Assignments"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Add(8, 2)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Assignments"

$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$tr8.Text = "Brown   -label computers. "
$null = $tr8.InsertAfter("`rMichael    -look into Machine Learning - what would be needed?")
$null = $tr8.InsertAfter("`rNick    -find feasibility of SQL tables vs. Pandas ")
$null = $tr8.InsertAfter("dataframe")
$null = $tr8.InsertAfter("`rPatti   -code editing/ documentation/ testing")

# ---------------------------------------------------------------------
# Step 6: new slide 9 -- "Data formats:" (the content that used to be on
# the original slide 4, now relocated to the very end of the deck).
# ---------------------------------------------------------------------
$s9 = $p.Slides.Add(9, 2)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Data formats:"

$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$tr9.Text = "Computer usage data:"
$null = $tr9.InsertAfter("`rCRR019,2016-07-01 14:10:24.793,in-use")
$null = $tr9.InsertAfter("`rWeather Data")
$null = $tr9.InsertAfter(":")
$null = $tr9.InsertAfter("`rFM-15 FEW:02 50 BKN:07 75 OVC:08 90 10 ")

$tr9.Paragraphs(2).IndentLevel = 2
$tr9.Paragraphs(4).IndentLevel = 2

# apply the "+mj-lt" (major latin / heading font) override to the raw
# data-sample runs, mirroring the source presentation.
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2).Font.Name = "+mj-lt"
$s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4).Font.Name = "+mj-lt"
$s9.Shapes.Item(2).TextFrame.TextRange.Runs(3, 1).Font.Name = "+mj-lt"
